# Scheduled-runner style refresh of the market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on a handful
# of leve rows across every job sheet. Values only -- no formulas, no
# structural changes.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 188.21053  # H55: 210 -> 188.21053
$ws.Cells.Item(55, 9).Value = 150.2  # I55: 183.33333 -> 150.2
$ws.Cells.Item(55, 10).Value = 201.78572  # J55: 216.66667 -> 201.78572
$ws.Cells.Item(55, 11).Value = 150.2  # K55: 183.33333 -> 150.2
$ws.Cells.Item(55, 12).Value = 201.78572  # L55: 216.66667 -> 201.78572
$ws.Cells.Item(55, 13).Value = 63.80000000000001  # M55: 30.66667000000001 -> 63.80000000000001
$ws.Cells.Item(55, 14).Value = -629.78572  # N55: -644.6666700000001 -> -629.78572
$ws.Cells.Item(74, 8).Value = 33338266  # H74: 25004700 -> 33338266
$ws.Cells.Item(74, 9).Value = 100000000  # I74: 50002000 -> 100000000
$ws.Cells.Item(74, 11).Value = 100000000  # K74: 50002000 -> 100000000
$ws.Cells.Item(74, 13).Value = -99999064  # M74: -50001064 -> -99999064
$ws.Cells.Item(77, 8).Value = 33338266  # H77: 25004700 -> 33338266
$ws.Cells.Item(77, 9).Value = 100000000  # I77: 50002000 -> 100000000
$ws.Cells.Item(77, 11).Value = 500000000  # K77: 250010000 -> 500000000
$ws.Cells.Item(77, 13).Value = -499995320  # M77: -250005320 -> -499995320
$ws.Cells.Item(112, 8).Value = 448480.84  # H112: 438315.38 -> 448480.84
$ws.Cells.Item(112, 10).Value = 482057.66  # J112: 470329.4 -> 482057.66
$ws.Cells.Item(112, 12).Value = 1446172.98  # L112: 1410988.2 -> 1446172.98
$ws.Cells.Item(112, 14).Value = -1448388.98  # N112: -1413204.2 -> -1448388.98
$ws.Cells.Item(113, 8).Value = 5269.643  # H113: 5537.923 -> 5269.643
$ws.Cells.Item(113, 9).Value = 3500  # I113: 3748 -> 3500
$ws.Cells.Item(113, 10).Value = 5405.769  # J113: 5863.364 -> 5405.769
$ws.Cells.Item(113, 11).Value = 3500  # K113: 3748 -> 3500
$ws.Cells.Item(113, 12).Value = 5405.769  # L113: 5863.364 -> 5405.769
$ws.Cells.Item(113, 13).Value = -246  # M113: -494 -> -246
$ws.Cells.Item(113, 14).Value = -11913.769  # N113: -12371.364 -> -11913.769
$ws.Cells.Item(125, 8).Value = 1220.25  # H125: 1254.125 -> 1220.25
$ws.Cells.Item(125, 10).Value = 1215.8889  # J125: 1266.6 -> 1215.8889
$ws.Cells.Item(125, 12).Value = 10943.0001  # L125: 11399.4 -> 10943.0001
$ws.Cells.Item(125, 14).Value = -15863.0001  # N125: -16319.4 -> -15863.0001
$ws.Cells.Item(129, 8).Value = 843.29895  # H129: 839.5599999999999 -> 843.29895
$ws.Cells.Item(129, 9).Value = 351.7  # I129: 350.54544 -> 351.7
$ws.Cells.Item(129, 10).Value = 899.8046000000001  # J129: 900 -> 899.8046000000001
$ws.Cells.Item(129, 11).Value = 1055.1  # K129: 1051.63632 -> 1055.1
$ws.Cells.Item(129, 12).Value = 2699.4138  # L129: 2700 -> 2699.4138
$ws.Cells.Item(129, 13).Value = 3944.9  # M129: 3948.36368 -> 3944.9
$ws.Cells.Item(129, 14).Value = -12699.4138  # N129: -12700 -> -12699.4138
$ws.Cells.Item(132, 8).Value = 375747.66  # H132: 921216.75 -> 375747.66
$ws.Cells.Item(132, 9).Value = 5781.0454  # I132: 14672.25 -> 5781.0454
$ws.Cells.Item(132, 10).Value = 2003600.8  # J132: 3338668.8 -> 2003600.8
$ws.Cells.Item(132, 11).Value = 17343.1362  # K132: 44016.75 -> 17343.1362
$ws.Cells.Item(132, 12).Value = 6010802.4  # L132: 10016006.4 -> 6010802.4
$ws.Cells.Item(132, 13).Value = -14813.1362  # M132: -41486.75 -> -14813.1362
$ws.Cells.Item(132, 14).Value = -6015862.4  # N132: -10021066.4 -> -6015862.4
$ws.Cells.Item(135, 8).Value = 446.58334  # H135: 428.57895 -> 446.58334
$ws.Cells.Item(135, 9).Value = 396.27274  # I135: 327.6875 -> 396.27274
$ws.Cells.Item(135, 10).Value = 1000  # J135: 966.6667 -> 1000
$ws.Cells.Item(135, 11).Value = 3566.45466  # K135: 2949.1875 -> 3566.45466
$ws.Cells.Item(135, 12).Value = 9000  # L135: 8700.0003 -> 9000
$ws.Cells.Item(135, 13).Value = -1031.45466  # M135: -414.1875 -> -1031.45466
$ws.Cells.Item(135, 14).Value = -14070  # N135: -13770.0003 -> -14070
$ws.Cells.Item(138, 8).Value = 3481.78  # H138: 3650.54 -> 3481.78
$ws.Cells.Item(138, 9).Value = 548  # I138: 521.3 -> 548
$ws.Cells.Item(138, 10).Value = 4799.855  # J138: 4991.643 -> 4799.855
$ws.Cells.Item(138, 11).Value = 1644  # K138: 1563.9 -> 1644
$ws.Cells.Item(138, 12).Value = 14399.565  # L138: 14974.929 -> 14399.565
$ws.Cells.Item(138, 13).Value = 3496  # M138: 3576.1 -> 3496
$ws.Cells.Item(138, 14).Value = -24679.565  # N138: -25254.929 -> -24679.565

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1002.63635  # H2: 1099.174 -> 1002.63635
$ws.Cells.Item(2, 9).Value = 989.25  # I2: 1130.1875 -> 989.25
$ws.Cells.Item(2, 10).Value = 1038.3334  # J2: 1028.2858 -> 1038.3334
$ws.Cells.Item(2, 11).Value = 989.25  # K2: 1130.1875 -> 989.25
$ws.Cells.Item(2, 12).Value = 1038.3334  # L2: 1028.2858 -> 1038.3334
$ws.Cells.Item(2, 13).Value = -876.25  # M2: -1017.1875 -> -876.25
$ws.Cells.Item(2, 14).Value = -1264.3334  # N2: -1254.2858 -> -1264.3334
$ws.Cells.Item(32, 8).Value = 3633.041  # H32: 3415.6538 -> 3633.041
$ws.Cells.Item(32, 9).Value = 3333.1  # I32: 3240.41 -> 3333.1
$ws.Cells.Item(32, 10).Value = 5017.385  # J32: 4044.4707 -> 5017.385
$ws.Cells.Item(32, 11).Value = 3333.1  # K32: 3240.41 -> 3333.1
$ws.Cells.Item(32, 12).Value = 5017.385  # L32: 4044.4707 -> 5017.385
$ws.Cells.Item(32, 13).Value = -3046.1  # M32: -2953.41 -> -3046.1
$ws.Cells.Item(32, 14).Value = -5591.385  # N32: -4618.4707 -> -5591.385
$ws.Cells.Item(74, 8).Value = 3316.0732  # H74: 3810.2646 -> 3316.0732
$ws.Cells.Item(74, 9).Value = 4400.905  # I74: 4576.95 -> 4400.905
$ws.Cells.Item(74, 10).Value = 2177  # J74: 2715 -> 2177
$ws.Cells.Item(74, 11).Value = 4400.905  # K74: 4576.95 -> 4400.905
$ws.Cells.Item(74, 12).Value = 2177  # L74: 2715 -> 2177
$ws.Cells.Item(74, 13).Value = -3526.905  # M74: -3702.95 -> -3526.905
$ws.Cells.Item(74, 14).Value = -3925  # N74: -4463 -> -3925
$ws.Cells.Item(77, 8).Value = 3316.0732  # H77: 3810.2646 -> 3316.0732
$ws.Cells.Item(77, 9).Value = 4400.905  # I77: 4576.95 -> 4400.905
$ws.Cells.Item(77, 10).Value = 2177  # J77: 2715 -> 2177
$ws.Cells.Item(77, 11).Value = 22004.525  # K77: 22884.75 -> 22004.525
$ws.Cells.Item(77, 12).Value = 10885  # L77: 13575 -> 10885
$ws.Cells.Item(77, 13).Value = -17636.525  # M77: -18516.75 -> -17636.525
$ws.Cells.Item(77, 14).Value = -19621  # N77: -22311 -> -19621
$ws.Cells.Item(116, 8).Value = 1002.63635  # H116: 1099.174 -> 1002.63635
$ws.Cells.Item(116, 9).Value = 989.25  # I116: 1130.1875 -> 989.25
$ws.Cells.Item(116, 10).Value = 1038.3334  # J116: 1028.2858 -> 1038.3334
$ws.Cells.Item(116, 11).Value = 989.25  # K116: 1130.1875 -> 989.25
$ws.Cells.Item(116, 12).Value = 1038.3334  # L116: 1028.2858 -> 1038.3334
$ws.Cells.Item(116, 13).Value = 1304.75  # M116: 1163.8125 -> 1304.75
$ws.Cells.Item(116, 14).Value = -5626.3334  # N116: -5616.2858 -> -5626.3334
$ws.Cells.Item(122, 8).Value = 1790.9642  # H122: 2015.9565 -> 1790.9642
$ws.Cells.Item(122, 9).Value = 1235.4286  # I122: 1347.75 -> 1235.4286
$ws.Cells.Item(122, 10).Value = 3457.5715  # J122: 3543.2856 -> 3457.5715
$ws.Cells.Item(122, 11).Value = 3706.2858  # K122: 4043.25 -> 3706.2858
$ws.Cells.Item(122, 12).Value = 10372.7145  # L122: 10629.8568 -> 10372.7145
$ws.Cells.Item(122, 13).Value = -1256.2858  # M122: -1593.25 -> -1256.2858
$ws.Cells.Item(122, 14).Value = -15272.7145  # N122: -15529.8568 -> -15272.7145
$ws.Cells.Item(132, 8).Value = 2295.353  # H132: 1815.7142 -> 2295.353
$ws.Cells.Item(132, 9).Value = 1234.1428  # I132: 1023.4722 -> 1234.1428
$ws.Cells.Item(132, 11).Value = 3702.4284  # K132: 3070.4166 -> 3702.4284
$ws.Cells.Item(132, 13).Value = -1172.4284  # M132: -540.4166 -> -1172.4284

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1002.63635  # H3: 1099.174 -> 1002.63635
$ws.Cells.Item(3, 9).Value = 989.25  # I3: 1130.1875 -> 989.25
$ws.Cells.Item(3, 10).Value = 1038.3334  # J3: 1028.2858 -> 1038.3334
$ws.Cells.Item(3, 11).Value = 989.25  # K3: 1130.1875 -> 989.25
$ws.Cells.Item(3, 12).Value = 1038.3334  # L3: 1028.2858 -> 1038.3334
$ws.Cells.Item(3, 13).Value = -875.25  # M3: -1016.1875 -> -875.25
$ws.Cells.Item(3, 14).Value = -1266.3334  # N3: -1256.2858 -> -1266.3334
$ws.Cells.Item(134, 8).Value = 2242.282  # H134: 2350.9714 -> 2242.282
$ws.Cells.Item(134, 9).Value = 1455  # I134: 1536.8 -> 1455
$ws.Cells.Item(134, 10).Value = 4013.6667  # J134: 4386.4 -> 4013.6667
$ws.Cells.Item(134, 11).Value = 4365  # K134: 4610.4 -> 4365
$ws.Cells.Item(134, 12).Value = 12041.0001  # L134: 13159.2 -> 12041.0001
$ws.Cells.Item(134, 13).Value = -1830  # M134: -2075.4 -> -1830
$ws.Cells.Item(134, 14).Value = -17111.0001  # N134: -18229.2 -> -17111.0001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2452.5667  # H31: 2639.037 -> 2452.5667
$ws.Cells.Item(31, 9).Value = 994.7778  # I31: 1022.125 -> 994.7778
$ws.Cells.Item(31, 10).Value = 4639.25  # J31: 4990.909 -> 4639.25
$ws.Cells.Item(31, 11).Value = 994.7778  # K31: 1022.125 -> 994.7778
$ws.Cells.Item(31, 12).Value = 4639.25  # L31: 4990.909 -> 4639.25
$ws.Cells.Item(31, 13).Value = -699.7778  # M31: -727.125 -> -699.7778
$ws.Cells.Item(31, 14).Value = -5229.25  # N31: -5580.909 -> -5229.25
$ws.Cells.Item(34, 8).Value = 2452.5667  # H34: 2639.037 -> 2452.5667
$ws.Cells.Item(34, 9).Value = 994.7778  # I34: 1022.125 -> 994.7778
$ws.Cells.Item(34, 10).Value = 4639.25  # J34: 4990.909 -> 4639.25
$ws.Cells.Item(34, 11).Value = 994.7778  # K34: 1022.125 -> 994.7778
$ws.Cells.Item(34, 12).Value = 4639.25  # L34: 4990.909 -> 4639.25
$ws.Cells.Item(34, 13).Value = -792.7778  # M34: -820.125 -> -792.7778
$ws.Cells.Item(34, 14).Value = -5043.25  # N34: -5394.909 -> -5043.25
$ws.Cells.Item(99, 8).Value = 11115070  # H99: 13337914 -> 11115070
$ws.Cells.Item(99, 9).Value = 22223584  # I99: 28573316 -> 22223584
$ws.Cells.Item(99, 10).Value = 6555.5557  # J99: 6937.5 -> 6555.5557
$ws.Cells.Item(99, 11).Value = 22223584  # K99: 28573316 -> 22223584
$ws.Cells.Item(99, 12).Value = 6555.5557  # L99: 6937.5 -> 6555.5557
$ws.Cells.Item(99, 13).Value = -22222086  # M99: -28571818 -> -22222086
$ws.Cells.Item(99, 14).Value = -9551.555700000001  # N99: -9933.5 -> -9551.555700000001
$ws.Cells.Item(107, 8).Value = 646.7826  # H107: 670.3182 -> 646.7826
$ws.Cells.Item(107, 9).Value = 540  # I107: 565.6875 -> 540
$ws.Cells.Item(107, 11).Value = 540  # K107: 565.6875 -> 540
$ws.Cells.Item(107, 13).Value = 1380  # M107: 1354.3125 -> 1380
$ws.Cells.Item(122, 8).Value = 2131.4  # H122: 2904.9 -> 2131.4
$ws.Cells.Item(122, 9).Value = 946.2857  # I122: 1300 -> 946.2857
$ws.Cells.Item(122, 10).Value = 3168.375  # J122: 3974.8333 -> 3168.375
$ws.Cells.Item(122, 11).Value = 2838.8571  # K122: 3900 -> 2838.8571
$ws.Cells.Item(122, 12).Value = 9505.125  # L122: 11924.4999 -> 9505.125
$ws.Cells.Item(122, 13).Value = -388.8571000000002  # M122: -1450 -> -388.8571000000002
$ws.Cells.Item(122, 14).Value = -14405.125  # N122: -16824.4999 -> -14405.125
$ws.Cells.Item(126, 8).Value = 11115070  # H126: 13337914 -> 11115070
$ws.Cells.Item(126, 9).Value = 22223584  # I126: 28573316 -> 22223584
$ws.Cells.Item(126, 10).Value = 6555.5557  # J126: 6937.5 -> 6555.5557
$ws.Cells.Item(126, 11).Value = 66670752  # K126: 85719948 -> 66670752
$ws.Cells.Item(126, 12).Value = 19666.6671  # L126: 20812.5 -> 19666.6671
$ws.Cells.Item(126, 13).Value = -66668282  # M126: -85717478 -> -66668282
$ws.Cells.Item(126, 14).Value = -24606.6671  # N126: -25752.5 -> -24606.6671
$ws.Cells.Item(132, 8).Value = 2583.9167  # H132: 2770.2727 -> 2583.9167
$ws.Cells.Item(132, 9).Value = 1655.85  # I132: 1780.5 -> 1655.85
$ws.Cells.Item(132, 11).Value = 4967.549999999999  # K132: 5341.5 -> 4967.549999999999
$ws.Cells.Item(132, 13).Value = -2437.549999999999  # M132: -2811.5 -> -2437.549999999999
$ws.Cells.Item(134, 8).Value = 9591.429  # H134: 6462.227 -> 9591.429
$ws.Cells.Item(134, 9).Value = 17372.666  # I134: 9028.166999999999 -> 17372.666
$ws.Cells.Item(134, 10).Value = 3755.5  # J134: 3383.1 -> 3755.5
$ws.Cells.Item(134, 11).Value = 52117.99800000001  # K134: 27084.501 -> 52117.99800000001
$ws.Cells.Item(134, 12).Value = 11266.5  # L134: 10149.3 -> 11266.5
$ws.Cells.Item(134, 13).Value = -49582.99800000001  # M134: -24549.501 -> -49582.99800000001
$ws.Cells.Item(134, 14).Value = -16336.5  # N134: -15219.3 -> -16336.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 771.86456  # H131: 5102848.5 -> 771.86456
$ws.Cells.Item(131, 9).Value = 334.66666  # I131: 125000250 -> 334.66666
$ws.Cells.Item(131, 10).Value = 785.9677  # J131: 831.34045 -> 785.9677
$ws.Cells.Item(131, 11).Value = 1003.99998  # K131: 375000750 -> 1003.99998
$ws.Cells.Item(131, 12).Value = 2357.9031  # L131: 2494.02135 -> 2357.9031
$ws.Cells.Item(131, 13).Value = 4036.00002  # M131: -374995710 -> 4036.00002
$ws.Cells.Item(131, 14).Value = -12437.9031  # N131: -12574.02135 -> -12437.9031

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2001.3334  # H132: 2226.5 -> 2001.3334
$ws.Cells.Item(132, 9).Value = 1492.6129  # I132: 1642.68 -> 1492.6129
$ws.Cells.Item(132, 10).Value = 2929  # J132: 3199.5334 -> 2929
$ws.Cells.Item(132, 11).Value = 4477.8387  # K132: 4928.04 -> 4477.8387
$ws.Cells.Item(132, 12).Value = 8787  # L132: 9598.600199999999 -> 8787
$ws.Cells.Item(132, 13).Value = -1947.8387  # M132: -2398.04 -> -1947.8387
$ws.Cells.Item(132, 14).Value = -13847  # N132: -14658.6002 -> -13847

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3441  # H7: 2838.6667 -> 3441
$ws.Cells.Item(7, 9).Value = 2085.9092  # I7: 1619.5 -> 2085.9092
$ws.Cells.Item(7, 10).Value = 5097.222  # J7: 6496.1665 -> 5097.222
$ws.Cells.Item(7, 11).Value = 2085.9092  # K7: 1619.5 -> 2085.9092
$ws.Cells.Item(7, 12).Value = 5097.222  # L7: 6496.1665 -> 5097.222
$ws.Cells.Item(7, 13).Value = -1973.9092  # M7: -1507.5 -> -1973.9092
$ws.Cells.Item(7, 14).Value = -5321.222  # N7: -6720.1665 -> -5321.222
$ws.Cells.Item(16, 8).Value = 484.875  # H16: 617.8182 -> 484.875
$ws.Cells.Item(16, 9).Value = 484.875  # I16: 617.8182 -> 484.875
$ws.Cells.Item(16, 11).Value = 484.875  # K16: 617.8182 -> 484.875
$ws.Cells.Item(16, 13).Value = -314.875  # M16: -447.8182 -> -314.875
$ws.Cells.Item(87, 8).Value = 50000  # H87: 45000 -> 50000
$ws.Cells.Item(87, 10).Value = 50000  # J87: 45000 -> 50000
$ws.Cells.Item(87, 12).Value = 50000  # L87: 45000 -> 50000
$ws.Cells.Item(87, 14).Value = -52246  # N87: -47246 -> -52246
$ws.Cells.Item(90, 8).Value = 50000  # H90: 45000 -> 50000
$ws.Cells.Item(90, 10).Value = 50000  # J90: 45000 -> 50000
$ws.Cells.Item(90, 12).Value = 150000  # L90: 135000 -> 150000
$ws.Cells.Item(90, 14).Value = -161232  # N90: -146232 -> -161232
$ws.Cells.Item(122, 8).Value = 2100.7896  # H122: 2579.6428 -> 2100.7896
$ws.Cells.Item(122, 9).Value = 1384.1666  # I122: 1624.2307 -> 1384.1666
$ws.Cells.Item(122, 11).Value = 4152.4998  # K122: 4872.6921 -> 4152.4998
$ws.Cells.Item(122, 13).Value = -1702.4998  # M122: -2422.6921 -> -1702.4998
$ws.Cells.Item(126, 8).Value = 3441  # H126: 2838.6667 -> 3441
$ws.Cells.Item(126, 9).Value = 2085.9092  # I126: 1619.5 -> 2085.9092
$ws.Cells.Item(126, 10).Value = 5097.222  # J126: 6496.1665 -> 5097.222
$ws.Cells.Item(126, 11).Value = 6257.7276  # K126: 4858.5 -> 6257.7276
$ws.Cells.Item(126, 12).Value = 15291.666  # L126: 19488.4995 -> 15291.666
$ws.Cells.Item(126, 13).Value = -3787.7276  # M126: -2388.5 -> -3787.7276
$ws.Cells.Item(126, 14).Value = -20231.666  # N126: -24428.4995 -> -20231.666
$ws.Cells.Item(132, 8).Value = 3687.4856  # H132: 4619.7036 -> 3687.4856
$ws.Cells.Item(132, 9).Value = 1330.95  # I132: 1629.2667 -> 1330.95
$ws.Cells.Item(132, 10).Value = 6829.533  # J132: 8357.75 -> 6829.533
$ws.Cells.Item(132, 11).Value = 3992.85  # K132: 4887.800099999999 -> 3992.85
$ws.Cells.Item(132, 12).Value = 20488.599  # L132: 25073.25 -> 20488.599
$ws.Cells.Item(132, 13).Value = -1462.85  # M132: -2357.800099999999 -> -1462.85
$ws.Cells.Item(132, 14).Value = -25548.599  # N132: -30133.25 -> -25548.599
$ws.Cells.Item(136, 8).Value = 2897.4878  # H136: 3188.1667 -> 2897.4878
$ws.Cells.Item(136, 9).Value = 1204.619  # I136: 1274.9412 -> 1204.619
$ws.Cells.Item(136, 10).Value = 4675  # J136: 4900 -> 4675
$ws.Cells.Item(136, 11).Value = 3613.857  # K136: 3824.8236 -> 3613.857
$ws.Cells.Item(136, 12).Value = 14025  # L136: 14700 -> 14025
$ws.Cells.Item(136, 13).Value = -1063.857  # M136: -1274.8236 -> -1063.857
$ws.Cells.Item(136, 14).Value = -19125  # N136: -19800 -> -19125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3723.7896  # H122: 3874.7646 -> 3723.7896
$ws.Cells.Item(122, 9).Value = 2125.1428  # I122: 2067 -> 2125.1428
$ws.Cells.Item(122, 10).Value = 8200  # J122: 9750 -> 8200
$ws.Cells.Item(122, 11).Value = 6375.428400000001  # K122: 6201 -> 6375.428400000001
$ws.Cells.Item(122, 12).Value = 24600  # L122: 29250 -> 24600
$ws.Cells.Item(122, 13).Value = -3925.428400000001  # M122: -3751 -> -3925.428400000001
$ws.Cells.Item(122, 14).Value = -29500  # N122: -34150 -> -29500
$ws.Cells.Item(126, 8).Value = 2296.9443  # H126: 2168.32 -> 2296.9443
$ws.Cells.Item(126, 9).Value = 1572.6923  # I126: 1508.5 -> 1572.6923
$ws.Cells.Item(126, 10).Value = 4180  # J126: 3008.0908 -> 4180
$ws.Cells.Item(126, 11).Value = 4718.0769  # K126: 4525.5 -> 4718.0769
$ws.Cells.Item(126, 12).Value = 12540  # L126: 9024.2724 -> 12540
$ws.Cells.Item(126, 13).Value = -2248.0769  # M126: -2055.5 -> -2248.0769
$ws.Cells.Item(126, 14).Value = -17480  # N126: -13964.2724 -> -17480
$ws.Cells.Item(132, 8).Value = 8334770.5  # H132: 11906545 -> 8334770.5
$ws.Cells.Item(132, 9).Value = 929  # I132: 1039.5217 -> 929
$ws.Cells.Item(132, 10).Value = 41670136  # J132: 66671868 -> 41670136
$ws.Cells.Item(132, 11).Value = 2787  # K132: 3118.5651 -> 2787
$ws.Cells.Item(132, 12).Value = 125010408  # L132: 200015604 -> 125010408
$ws.Cells.Item(132, 13).Value = -588.5650999999998  # M132: -257 -> -588.5650999999998
$ws.Cells.Item(132, 14).Value = -125015468  # N132: -200020664 -> -125015468
